$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_QueryLog_External")

# --- Fill in hand-typed Effort (PDs) / Comments values that replace the
# --- formula-derived 0s on a handful of rows (author overrode the computed
# --- effort with a literal estimate + a clarifying comment). ---
$ws.Range("N11").Value = "24"
$ws.Range("O11").Value = "Assumptions : RP will expose 2 apis`n1. to provide applicant info.`n2. receive updated info and incorporate it."

$ws.Range("N12").Value = "36"
$ws.Range("O12").Value = "Estimation may change after understanding overall scope of the change."

$ws.Range("N13").Value = "52"
$ws.Range("O13").Value = "Need more clarification on the requirement. Estimation may change after clarification."

$ws.Range("N14").Value = "30"

$ws.Range("N15").Value = "28"
$ws.Range("O15").Value = "Estimation may change after understanding overall scope of the change."

$ws.Range("N20").Value = "20"
$ws.Range("O20").Value = "Since the requirement is not detailed the effort may change."

$ws.Range("N40").Value = "45"

$ws.Range("N41").Value = "55"
$ws.Range("O41").Value = "Need more clarification on the requirement. Estimation may change after clarification."

$ws.Range("N42").Value = "12"
$ws.Range("O42").Value = "Change algorithm from lavenstine distance to phonetic and soundex match"

# --- Apply an AutoFilter on the full data range (rows grew from 34 to 53
# --- data rows), restricting column D ("Module") to "Registration Processor"
# --- only -- this hides every row whose Module isn't Registration Processor. ---
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$filterRange = $ws.Range("A2:H53")
$filterRange.AutoFilter(4, @("Registration Processor"), 7)

# Keep the _FilterDatabase defined name's range in sync with the new filter
# extent (Excel does this automatically; make it explicit here too).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=MOSIP_QueryLog_External!`$A`$2:`$H`$53"
    }
}
